# Applies the "login_test_data.xlsx" update:
#  - adds 3 new worksheets (homePageData, listPurchaseData, userCreation)
#  - adds a welcome banner row to loginData
#  - re-colors / re-styles a few existing cells on loginData
#  - populates the new sheets with their data + a mailto hyperlink

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("loginData")

# ---------------------------------------------------------------------------
# 1. loginData: restyle rows 1-4, add new "Welcome admin," row 5
# ---------------------------------------------------------------------------

# bump the font used by the two highlighted rows from 10pt -> 11pt
$ws1.Range("A1:B2").Font.Size = 11

# criss-cross the existing highlight colors onto rows 3-4 (A3/B4 pick up the
# "admin" green, B3/A4 pick up the "saranya" orange)
$ws1.Range("A3").Interior.Color = $ws1.Range("A1").Interior.Color
$ws1.Range("B4").Interior.Color = $ws1.Range("A1").Interior.Color
$ws1.Range("B3").Interior.Color = $ws1.Range("A2").Interior.Color
$ws1.Range("A4").Interior.Color = $ws1.Range("A2").Interior.Color

# uniform row height bump for the original 4 rows
$ws1.Rows.Item(1).RowHeight = 13.8
$ws1.Rows.Item(2).RowHeight = 13.8
$ws1.Rows.Item(3).RowHeight = 13.8
$ws1.Rows.Item(4).RowHeight = 13.8

# new banner row
$ws1.Range("A5").Value = "Welcome admin,"
$ws1.Range("A5").Font.Name = "Monospace"
$ws1.Range("A5").Font.Size = 12
$ws1.Range("A5").Font.Color = 16711722  # BGR for FF2A00FF
$ws1.Rows.Item(5).RowHeight = 15

$excel.ActiveWindow.Zoom = 90
$ws1.Range("E12").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. homePageData
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $last)
$ws2.Name = "homePageData"

$ws2.Range("A1").Value = "Maxfest Enterprises PVT LTD,Kochi"
$ws2.Range("A2").Value = "Calculator"
$ws2.Range("A1:A2").Font.Name = "Monospace"
$ws2.Range("A1:A2").Font.Size = 12
$ws2.Range("A1:A2").Font.Color = 16711722
$ws2.Rows.Item(1).RowHeight = 15
$ws2.Rows.Item(2).RowHeight = 15
$ws2.Columns.Item(1).ColumnWidth = 51.82
$ws2.Range("A13").Select() | Out-Null

# ---------------------------------------------------------------------------
# 3. listPurchaseData
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $last)
$ws3.Name = "listPurchaseData"

$ws3.Range("A1").Value = "Purchases"
$ws3.Range("A1").Font.Name = "Monospace"
$ws3.Range("A1").Font.Size = 12
$ws3.Range("A1").Font.Color = 16711722
$ws3.Rows.Item(1).RowHeight = 15
$ws3.Columns.Item(1).ColumnWidth = 15.8
$excel.ActiveWindow.Zoom = 90
$ws3.Range("C23").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. userCreation
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $last)
$ws4.Name = "userCreation"

$ws4.Range("A1").Value = "Diya"
$ws4.Range("B1").Value = "Diyaninnan+1@gmail.com"
$ws4.Range("C1").Value = "diya231"
$ws4.Range("D1").Value = "diya231"

# Add the hyperlink BEFORE the font styling below - Excel re-paints a cell
# with its built-in "Hyperlink" style as soon as a hyperlink lands on it,
# which would otherwise clobber the Monospace/colored font we want here.
$ws4.Hyperlinks.Add($ws4.Range("B1"), "mailto:Diyaninnan+1@gmail.com", "", "", "Diyaninnan+1@gmail.com") | Out-Null

$ws4.Range("A1:D1").Font.Name = "Monospace"
$ws4.Range("A1:D1").Font.Size = 12
$ws4.Range("A1:D1").Font.Underline = 0
$ws4.Range("A1:D1").Font.Color = 16711722
$ws4.Rows.Item(1).RowHeight = 15
$ws4.Columns.Item(2).ColumnWidth = 30.07

$excel.ActiveWindow.Zoom = 90
$ws4.Range("J19").Select() | Out-Null

# ---------------------------------------------------------------------------
# Re-select loginData as the active sheet/tab, matching the source workbook.
# ---------------------------------------------------------------------------
$ws1.Activate()
